$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.01486144694607106
$ws.Range("D2").Value = 0.02060685928102401
$ws.Range("E2").Value = 0.094008856121512
$ws.Range("F2").Value = 0.7601367278425784
$ws.Range("G2").Value = 0.6052006538951105
$ws.Range("H2").Value = 0.6993633722407964
$ws.Range("I2").Value = 0.5712109158885781
$ws.Range("K2").Value = 1.195645398134332
$ws.Range("M2").Value = 0.3868638408174476

# Row 3
$ws.Range("C3").Value = 0.01327159670282185
$ws.Range("D3").Value = 0.01998618411921171
$ws.Range("E3").Value = 0.08978452448139507
$ws.Range("F3").Value = 0.7547478513771892
$ws.Range("G3").Value = 0.6008226339211831
$ws.Range("H3").Value = 0.7035600338712555
$ws.Range("I3").Value = 0.5680939701469683
$ws.Range("K3").Value = 1.049026711075612
$ws.Range("M3").Value = 0.346003700589101

# Row 4
$ws.Range("C4").Value = 0.01228988104438855
$ws.Range("D4").Value = 0.01960754166852752
$ws.Range("E4").Value = 0.08729285887928739
$ws.Range("F4").Value = 0.7521823674383796
$ws.Range("G4").Value = 0.5988333556774421
$ws.Range("H4").Value = 0.7066747388953161
$ws.Range("I4").Value = 0.5667457098550557
$ws.Range("K4").Value = 0.9588623443301287
$ws.Range("M4").Value = 0.321002343968523

# Row 5
$ws.Range("C5").Value = 0.01188845162072028
$ws.Range("D5").Value = 0.01945387516252595
$ws.Range("E5").Value = 0.08630284111840325
$ws.Range("F5").Value = 0.7513229264827999
$ws.Range("G5").Value = 0.5981973286668989
$ws.Range("H5").Value = 0.7080788845903641
$ws.Range("I5").Value = 0.5663378018455063
$ws.Range("K5").Value = 0.9220855925075568
$ws.Range("M5").Value = 0.3108357594355695

# Row 6
$ws.Range("C6").Value = 0.01182171231657492
$ws.Range("D6").Value = 0.01942839769709082
$ws.Range("E6").Value = 0.08613997167145016
$ws.Range("F6").Value = 0.7511914255595684
$ws.Range("G6").Value = 0.5981022309884025
$ws.Range("H6").Value = 0.708320176791446
$ws.Range("I6").Value = 0.5662785960930989
$ws.Range("K6").Value = 0.9159768192673141
$ws.Range("M6").Value = 0.3091489094102897

# Row 7
$ws.Range("C7").Value = 0.01228447275506994
$ws.Range("D7").Value = 0.01960546668218299
$ws.Range("E7").Value = 0.08727940490312136
$ws.Range("F7").Value = 0.7521700247437195
$ws.Range("G7").Value = 0.5988240724394984
$ws.Range("H7").Value = 0.7066931301021384
$ws.Range("I7").Value = 0.5667396365845292
$ws.Range("K7").Value = 0.9583664964806076
$ws.Range("M7").Value = 0.3208651463484316

# Row 8
$ws.Range("C8").Value = 0.01431442867765043
$ws.Range("D8").Value = 0.02039235029437592
$ws.Range("E8").Value = 0.09253092600503976
$ws.Range("F8").Value = 0.7581238263123993
$ws.Range("G8").Value = 0.6035453966349849
$ws.Range("H8").Value = 0.7006984959451898
$ws.Range("I8").Value = 0.5700183907189995
$ws.Range("K8").Value = 1.145120939149137
$ws.Range("M8").Value = 0.3727570921326659

# Row 9
$ws.Range("C9").Value = 0.01825047879609087
$ws.Range("D9").Value = 0.02195430972687262
$ws.Range("E9").Value = 0.103654037010628
$ws.Range("F9").Value = 0.7757392793494375
$ws.Range("G9").Value = 0.6184003848310908
$ws.Range("H9").Value = 0.6932295000653426
$ws.Range("I9").Value = 0.5809680616964314
$ws.Range("K9").Value = 1.510205360288921
$ws.Range("M9").Value = 0.4752210225190083

# Row 10
$ws.Range("C10").Value = 0.02111438516911335
$ws.Range("D10").Value = 0.02311276618208069
$ws.Range("E10").Value = 0.1123505117302344
$ws.Range("F10").Value = 0.7923636989575726
$ws.Range("G10").Value = 0.6327994036293063
$ws.Range("H10").Value = 0.6903811143289573
$ws.Range("I10").Value = 0.5918152763761384
$ws.Range("K10").Value = 1.777725686464123
$ws.Range("M10").Value = 0.5509582563475277

# Row 11
$ws.Range("C11").Value = 0.0224110746141406
$ws.Range("D11").Value = 0.02364203087070393
$ws.Range("E11").Value = 0.1164253115327654
$ws.Range("F11").Value = 0.8007394899827602
$ws.Range("G11").Value = 0.6401225384238387
$ws.Range("H11").Value = 0.689664021585358
$ws.Range("I11").Value = 0.5973687149111271
$ws.Range("K11").Value = 1.899273907596296
$ws.Range("M11").Value = 0.5855189344148783

# Row 12
$ws.Range("C12").Value = 0.02290120243426941
$ws.Range("D12").Value = 0.02384276613185676
$ws.Range("E12").Value = 0.1179857706730871
$ws.Range("F12").Value = 0.8040291574257452
$ws.Range("G12").Value = 0.6430080154789835
$ws.Range("H12").Value = 0.6894761372888496
$ws.Range("I12").Value = 0.5995614604278856
$ws.Range("K12").Value = 1.945279178780083
$ws.Range("M12").Value = 0.5986219904620498

# Row 13
$ws.Range("C13").Value = 0.02279568500804174
$ws.Range("D13").Value = 0.02379952042996791
$ws.Range("E13").Value = 0.1176489175167035
$ws.Range("F13").Value = 0.8033154075877178
$ws.Range("G13").Value = 0.642381559531529
$ws.Range("H13").Value = 0.6895128733492157
$ws.Range("I13").Value = 0.5990852082403251
$ws.Range("K13").Value = 1.935372133995088
$ws.Range("M13").Value = 0.5957993118941118

# Row 14
$ws.Range("C14").Value = 0.02245141588686295
$ws.Range("D14").Value = 0.02365853925718397
$ws.Range("E14").Value = 0.1165533403046624
$ws.Range("F14").Value = 0.8010077635817936
$ws.Range("G14").Value = 0.6403576694288091
$ws.Range("H14").Value = 0.6896468848929516
$ws.Range("I14").Value = 0.5975473096901922
$ws.Range("K14").Value = 1.903059244101485
$ws.Range("M14").Value = 0.5865966151777116

# Row 15
$ws.Range("C15").Value = 0.0222404232476805
$ws.Range("D15").Value = 0.02357222475782095
$ws.Range("E15").Value = 0.1158845467021692
$ws.Range("F15").Value = 0.79960965335043
$ws.Range("G15").Value = 0.6391326476333745
$ws.Range("H15").Value = 0.6897398797382124
$ws.Range("I15").Value = 0.5966170171400194
$ws.Range("K15").Value = 1.883263716507088
$ws.Range("M15").Value = 0.5809617463330312

# Row 16
$ws.Range("C16").Value = 0.0210295186009759
$ws.Range("D16").Value = 0.02307822216306477
$ws.Range("E16").Value = 0.112086632795048
$ws.Range("F16").Value = 0.791832770886657
$ws.Range("G16").Value = 0.6323364771406119
$ws.Range("H16").Value = 0.690439660732352
$ws.Range("I16").Value = 0.5914648674317746
$ws.Range("K16").Value = 1.76977914869309
$ws.Range("M16").Value = 0.5487018158785446

# Row 17
$ws.Range("C17").Value = 0.02028508638039028
$ws.Range("D17").Value = 0.02277574113077918
$ws.Range("E17").Value = 0.1097873918263161
$ws.Range("F17").Value = 0.787270905506503
$ws.Range("G17").Value = 0.6283660818979229
$ws.Range("H17").Value = 0.6910174827021933
$ws.Range("I17").Value = 0.5884632746879674
$ws.Range("K17").Value = 1.700121308789562
$ws.Range("M17").Value = 0.5289390687483717

# Row 18
$ws.Range("C18").Value = 0.01985633421045918
$ws.Range("D18").Value = 0.02260197711138545
$ws.Range("E18").Value = 0.108476083720852
$ws.Range("F18").Value = 0.7847235296560626
$ws.Range("G18").Value = 0.6261550964249949
$ws.Range("H18").Value = 0.6914042689002429
$ws.Range("I18").Value = 0.5867950537820832
$ws.Range("K18").Value = 1.660042083692417
$ws.Range("M18").Value = 0.5175821568626304

# Row 19
$ws.Range("C19").Value = 0.01971106807570777
$ws.Range("D19").Value = 0.02254318095287289
$ws.Range("E19").Value = 0.1080340025661073
$ws.Range("F19").Value = 0.783874138737616
$ws.Range("G19").Value = 0.6254189407185464
$ws.Range("H19").Value = 0.6915445640847224
$ws.Range("I19").Value = 0.5862401980477827
$ws.Range("K19").Value = 1.646469597423902
$ws.Range("M19").Value = 0.5137386290910797

# Row 20
$ws.Range("C20").Value = 0.02036439206598573
$ws.Range("D20").Value = 0.02280791859276121
$ws.Range("E20").Value = 0.1100309929934369
$ws.Range("F20").Value = 0.7877486004388885
$ws.Range("G20").Value = 0.6287812062705171
$ws.Range("H20").Value = 0.6909503351470363
$ws.Range("I20").Value = 0.5887767686840704
$ws.Range("K20").Value = 1.707537949169989
$ws.Range("M20").Value = 0.5310418002233206

# Row 21
$ws.Range("C21").Value = 0.02255256067099509
$ws.Range("D21").Value = 0.0236999404198599
$ws.Range("E21").Value = 0.1168746623084118
$ws.Range("F21").Value = 0.8016823658733614
$ws.Range("G21").Value = 0.6409490757641834
$ws.Range("H21").Value = 0.689605248235452
$ws.Range("I21").Value = 0.5979965850617148
$ws.Range("K21").Value = 1.912550939732398
$ws.Range("M21").Value = 0.5892992429279502

# Row 22
$ws.Range("C22").Value = 0.02397739820273159
$ws.Range("D22").Value = 0.02428475408758857
$ws.Range("E22").Value = 0.1214491011784062
$ws.Range("F22").Value = 0.8114767764948851
$ws.Range("G22").Value = 0.649556995710185
$ws.Range("H22").Value = 0.689214020868846
$ws.Range("I22").Value = 0.6045459261200961
$ws.Range("K22").Value = 2.046407718239095
$ws.Range("M22").Value = 0.6274653422941157

# Row 23
$ws.Range("C23").Value = 0.02321742339168509
$ws.Range("D23").Value = 0.02397246514286877
$ws.Range("E23").Value = 0.1189982155074958
$ws.Range("F23").Value = 0.8061860452389027
$ws.Range("G23").Value = 0.6449024051044319
$ws.Range("H23").Value = 0.6893780362994875
$ws.Range("I23").Value = 0.6010022501941492
$ws.Range("K23").Value = 1.974978183178564
$ws.Range("M23").Value = 0.6070869269314443

# Row 24
$ws.Range("C24").Value = 0.02032854035623899
$ws.Range("D24").Value = 0.02279337073570531
$ws.Range("E24").Value = 0.1099208280525872
$ws.Range("F24").Value = 0.7875324000799964
$ws.Range("G24").Value = 0.6285933054599155
$ws.Range("H24").Value = 0.6909805225509444
$ws.Range("I24").Value = 0.5886348592039354
$ws.Range("K24").Value = 1.704184985668178
$ws.Range("M24").Value = 0.5300911399122725

# Row 25
$ws.Range("C25").Value = 0.01719052790931386
$ws.Range("D25").Value = 0.02152980433108809
$ws.Range("E25").Value = 0.1005544209565841
$ws.Range("F25").Value = 0.7703315938794333
$ws.Range("G25").Value = 0.6137750031266194
$ws.Range("H25").Value = 0.6947884894249796
$ws.Range("I25").Value = 0.5775171342388674
$ws.Range("K25").Value = 1.411563825335122
$ws.Range("M25").Value = 0.4474236091305812

